$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F3").Value = 1.67
$ws.Range("G3").Value = 1.85
$ws.Range("H3").Value = 2.18
$ws.Range("I3").Value = 11.5
$ws.Range("J3").Value = 3.2
$ws.Range("L3").Value = 1.37
$ws.Range("P3").Value = 1.62
$ws.Range("Q3").Value = 2.06
$ws.Range("S3").Value = 3.6
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 13
$ws.Range("AC4").Value = 15
$ws.Range("AD4").Value = 38
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 12
$ws.Range("AG4").Value = 12.5
$ws.Range("AH4").Value = 980
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 16
$ws.Range("AK4").Value = 17
$ws.Range("AL4").Value = 32
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 5.6
$ws.Range("AO4").Value = 1000
$ws.Range("F4").Value = 1.46
$ws.Range("G4").Value = 1.48
$ws.Range("I4").Value = 9.199999999999999
$ws.Range("J4").Value = 4.9
$ws.Range("L4").Value = 1.25
$ws.Range("N4").Value = 4.6
$ws.Range("P4").Value = 2.48
$ws.Range("Q4").Value = 1.54
$ws.Range("R4").Value = 1.59
$ws.Range("S4").Value = 2.16
$ws.Range("U4").Value = 2.04
$ws.Range("V4").Value = 1.12
$ws.Range("W4").Value = 3.05
$ws.Range("X4").Value = 30
$ws.Range("Y4").Value = 42
$ws.Range("Z4").Value = 1000
$ws.Range("F5").Value = 2.94
$ws.Range("G5").Value = 1000
$ws.Range("H5").Value = 1.83
$ws.Range("I5").Value = 2.26
$ws.Range("J5").Value = 3.6
$ws.Range("P5").Value = 1.84
$ws.Range("Q5").Value = 1.7
$ws.Range("R5").Value = 1.1
$ws.Range("S5").Value = 1.05
$ws.Range("V5").Value = 1.79
$ws.Range("F6").Value = 1.14
$ws.Range("H6").Value = 1.44
$ws.Range("N6").Value = 1.11
$ws.Range("P6").Value = 1.25
$ws.Range("S6").Value = 1.19
$ws.Range("O8").Value = 1.21
$ws.Range("Q8").Value = 1.21
$ws.Range("S8").Value = 1.21
$ws.Range("O10").Value = 1.15
$ws.Range("Q10").Value = 1.15
$ws.Range("S10").Value = 1.15
$ws.Range("H11").Value = 5.1
$ws.Range("J11").Value = 3.1
$ws.Range("P11").Value = 1.48
$ws.Range("Q11").Value = 2.12
$ws.Range("S11").Value = 2.12
$ws.Range("AC12").Value = 8.800000000000001
$ws.Range("AF12").Value = 980
$ws.Range("F12").Value = 2.84
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 2.82
$ws.Range("I12").Value = 3.1
$ws.Range("J12").Value = 3
$ws.Range("N12").Value = 2.68
$ws.Range("O12").Value = 1.51
$ws.Range("P12").Value = 1.58
$ws.Range("R12").Value = 1.2
$ws.Range("T12").Value = 2.04
$ws.Range("U12").Value = 1.83
$ws.Range("W12").Value = 1.5
$ws.Range("AC13").Value = 980
$ws.Range("AD13").Value = 22
$ws.Range("AE13").Value = 1000
$ws.Range("AF13").Value = 11.5
$ws.Range("AG13").Value = 980
$ws.Range("AH13").Value = 980
$ws.Range("AJ13").Value = 980
$ws.Range("AK13").Value = 980
$ws.Range("AN13").Value = 23
$ws.Range("F13").Value = 1.94
$ws.Range("G13").Value = 2.04
$ws.Range("I13").Value = 5
$ws.Range("J13").Value = 3.4
$ws.Range("T13").Value = 1.83
$ws.Range("U13").Value = 1.77
$ws.Range("V13").Value = 1.25
$ws.Range("W13").Value = 2
$ws.Range("F14").Value = 2.2
$ws.Range("G14").Value = 2.4
$ws.Range("H14").Value = 3.8
$ws.Range("I14").Value = 4.4
$ws.Range("J14").Value = 2.92
$ws.Range("K14").Value = 3.55
$ws.Range("O14").Value = 1.51
$ws.Range("P14").Value = 1.57
$ws.Range("Q14").Value = 2.32
$ws.Range("S14").Value = 5.2
$ws.Range("V14").Value = 1.3
$ws.Range("W14").Value = 1.71
$ws.Range("AB15").Value = 9.199999999999999
$ws.Range("AC15").Value = 8.4
$ws.Range("AD15").Value = 17
$ws.Range("AF15").Value = 16.5
$ws.Range("AG15").Value = 13.5
$ws.Range("AH15").Value = 29
$ws.Range("AI15").Value = 100
$ws.Range("AK15").Value = 38
$ws.Range("AM15").Value = 230
$ws.Range("AN15").Value = 42
$ws.Range("F15").Value = 2.42
$ws.Range("G15").Value = 2.6
$ws.Range("H15").Value = 3.4
$ws.Range("I15").Value = 3.9
$ws.Range("J15").Value = 2.94
$ws.Range("K15").Value = 3.35
$ws.Range("L15").Value = 1.59
$ws.Range("M15").Value = 1.13
$ws.Range("N15").Value = 2.6
$ws.Range("O15").Value = 1.53
$ws.Range("P15").Value = 1.53
$ws.Range("Q15").Value = 2.6
$ws.Range("R15").Value = 1.19
$ws.Range("S15").Value = 5.4
$ws.Range("T15").Value = 2.1
$ws.Range("U15").Value = 1.74
$ws.Range("V15").Value = 1.35
$ws.Range("W15").Value = 1.62
$ws.Range("X15").Value = 8.6
$ws.Range("Y15").Value = 10.5
$ws.Range("Z15").Value = 24
